# chore: update Sheets via scheduled runner
# Refresh the market-price-derived columns (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) on each
# per-job leve sheet with the latest pulled values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 43478660
$ws.Range("J33").Value = 96.666664
$ws.Range("L33").Value = 96.666664
$ws.Range("N33").Value = -554.666664
$ws.Range("H40").Value = 2172.7273
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2211.111
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2211.111
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2561.111
$ws.Range("H62").Value = 2150.5715
$ws.Range("I62").Value = 3001.25
$ws.Range("J62").Value = 1016.3333
$ws.Range("K62").Value = 3001.25
$ws.Range("L62").Value = 1016.3333
$ws.Range("M62").Value = -2377.25
$ws.Range("N62").Value = -2264.3333
$ws.Range("H64").Value = 27780428
$ws.Range("I64").Value = 55557356
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 55557356
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -55557108
$ws.Range("N64").Value = -3996
$ws.Range("H65").Value = 2150.5715
$ws.Range("I65").Value = 3001.25
$ws.Range("J65").Value = 1016.3333
$ws.Range("K65").Value = 15006.25
$ws.Range("L65").Value = 5081.6665
$ws.Range("M65").Value = -11886.25
$ws.Range("N65").Value = -11321.6665
$ws.Range("H67").Value = 27780428
$ws.Range("I67").Value = 55557356
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 55557356
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -55556498
$ws.Range("N67").Value = -5216
$ws.Range("H106").Value = 1003090.5
$ws.Range("I106").Value = 558989.4399999999
$ws.Range("K106").Value = 558989.4399999999
$ws.Range("M106").Value = -558358.4399999999
$ws.Range("H125").Value = 1625.75
$ws.Range("I125").Value = 2137.4546
$ws.Range("J125").Value = 500
$ws.Range("K125").Value = 19237.0914
$ws.Range("L125").Value = 4500
$ws.Range("M125").Value = -16777.0914
$ws.Range("N125").Value = -9420
$ws.Range("H138").Value = 8069218
$ws.Range("I138").Value = 1960.7368
$ws.Range("J138").Value = 20842376
$ws.Range("K138").Value = 5882.2104
$ws.Range("L138").Value = 62527128
$ws.Range("M138").Value = -742.2103999999999
$ws.Range("N138").Value = -62537408

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10715.213
$ws.Range("I32").Value = 12431.629
$ws.Range("J32").Value = 5709
$ws.Range("K32").Value = 12431.629
$ws.Range("L32").Value = 5709
$ws.Range("M32").Value = -12144.629
$ws.Range("N32").Value = -6283
$ws.Range("H122").Value = 4656.8486
$ws.Range("I122").Value = 6121.909
$ws.Range("J122").Value = 1726.7273
$ws.Range("K122").Value = 18365.727
$ws.Range("L122").Value = 5180.1819
$ws.Range("M122").Value = -15915.727
$ws.Range("N122").Value = -10080.1819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 846
$ws.Range("I64").Value = 1341.5
$ws.Range("K64").Value = 1341.5
$ws.Range("M64").Value = -1116.5
$ws.Range("H67").Value = 846
$ws.Range("I67").Value = 1341.5
$ws.Range("K67").Value = 1341.5
$ws.Range("M67").Value = -561.5
$ws.Range("H75").Value = 25120.4
$ws.Range("I75").Value = 19214.285
$ws.Range("J75").Value = 30288.25
$ws.Range("K75").Value = 19214.285
$ws.Range("L75").Value = 30288.25
$ws.Range("M75").Value = -18278.285
$ws.Range("N75").Value = -32160.25
$ws.Range("H78").Value = 25120.4
$ws.Range("I78").Value = 19214.285
$ws.Range("J78").Value = 30288.25
$ws.Range("K78").Value = 57642.855
$ws.Range("L78").Value = 90864.75
$ws.Range("M78").Value = -52962.855
$ws.Range("N78").Value = -100224.75
$ws.Range("H112").Value = 38635.668
$ws.Range("J112").Value = 38635.668
$ws.Range("L112").Value = 38635.668
$ws.Range("N112").Value = -41589.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2384.2856
$ws.Range("H65").Value = 2384.2856
$ws.Range("H70").Value = 33835
$ws.Range("J70").Value = 33835
$ws.Range("L70").Value = 33835
$ws.Range("N70").Value = -34465
$ws.Range("H73").Value = 33835
$ws.Range("J73").Value = 33835
$ws.Range("L73").Value = 33835
$ws.Range("N73").Value = -36019
$ws.Range("H99").Value = 1043.3793
$ws.Range("I99").Value = 927.1667
$ws.Range("J99").Value = 1233.5454
$ws.Range("K99").Value = 927.1667
$ws.Range("L99").Value = 1233.5454
$ws.Range("M99").Value = 570.8333
$ws.Range("N99").Value = -4229.5454
$ws.Range("H126").Value = 1043.3793
$ws.Range("I126").Value = 927.1667
$ws.Range("J126").Value = 1233.5454
$ws.Range("K126").Value = 2781.5001
$ws.Range("L126").Value = 3700.6362
$ws.Range("M126").Value = -311.5001000000002
$ws.Range("N126").Value = -8640.636200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 832.383
$ws.Range("I68").Value = 557.0135
$ws.Range("J68").Value = 1851.25
$ws.Range("K68").Value = 1671.0405
$ws.Range("L68").Value = 5553.75
$ws.Range("M68").Value = -860.0405000000001
$ws.Range("N68").Value = -7175.75
$ws.Range("H71").Value = 832.383
$ws.Range("I71").Value = 557.0135
$ws.Range("J71").Value = 1851.25
$ws.Range("K71").Value = 5013.1215
$ws.Range("L71").Value = 16661.25
$ws.Range("M71").Value = -957.1215000000002
$ws.Range("N71").Value = -24773.25
$ws.Range("H122").Value = 1579.1
$ws.Range("I122").Value = 2267.6
$ws.Range("J122").Value = 890.6
$ws.Range("K122").Value = 20408.4
$ws.Range("L122").Value = 8015.400000000001
$ws.Range("M122").Value = -17958.4
$ws.Range("N122").Value = -12915.4
$ws.Range("H132").Value = 1920.9706
$ws.Range("I132").Value = 840.5333000000001
$ws.Range("J132").Value = 2773.9473
$ws.Range("K132").Value = 7564.7997
$ws.Range("L132").Value = 24965.5257
$ws.Range("M132").Value = -5034.7997
$ws.Range("N132").Value = -30025.5257

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11113252
$ws.Range("I122").Value = 33333808
$ws.Range("J122").Value = 2975
$ws.Range("K122").Value = 100001424
$ws.Range("L122").Value = 8925
$ws.Range("M122").Value = -99998974
$ws.Range("N122").Value = -13825

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5080.1143
$ws.Range("I7").Value = 4780
$ws.Range("J7").Value = 5480.2666
$ws.Range("K7").Value = 4780
$ws.Range("L7").Value = 5480.2666
$ws.Range("M7").Value = -4668
$ws.Range("N7").Value = -5704.2666
$ws.Range("H68").Value = 2056.6667
$ws.Range("I68").Value = 1580
$ws.Range("J68").Value = 2152
$ws.Range("K68").Value = 1580
$ws.Range("L68").Value = 2152
$ws.Range("M68").Value = -831
$ws.Range("N68").Value = -3650
$ws.Range("H71").Value = 2056.6667
$ws.Range("I71").Value = 1580
$ws.Range("J71").Value = 2152
$ws.Range("K71").Value = 7900
$ws.Range("L71").Value = 10760
$ws.Range("M71").Value = -4156
$ws.Range("N71").Value = -18248
$ws.Range("H126").Value = 5080.1143
$ws.Range("I126").Value = 4780
$ws.Range("J126").Value = 5480.2666
$ws.Range("K126").Value = 14340
$ws.Range("L126").Value = 16440.7998
$ws.Range("M126").Value = -11870
$ws.Range("N126").Value = -21380.7998
$ws.Range("H132").Value = 9811040
$ws.Range("I132").Value = 5287.4614
$ws.Range("J132").Value = 20009024
$ws.Range("K132").Value = 15862.3842
$ws.Range("L132").Value = 60027072
$ws.Range("M132").Value = -13332.3842
$ws.Range("N132").Value = -60032132

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3265.375
$ws.Range("I96").Value = 2345
$ws.Range("K96").Value = 2345
$ws.Range("M96").Value = -972
$ws.Range("H113").Value = 2185
$ws.Range("I113").Value = 2366.4
$ws.Range("J113").Value = 2033.8334
$ws.Range("K113").Value = 7099.200000000001
$ws.Range("L113").Value = 6101.5002
$ws.Range("M113").Value = -4929.200000000001
$ws.Range("N113").Value = -10441.5002
$ws.Range("H126").Value = 3506.3
$ws.Range("I126").Value = 1889.7646
$ws.Range("J126").Value = 12666.667
$ws.Range("K126").Value = 5669.293799999999
$ws.Range("L126").Value = 38000.001
$ws.Range("M126").Value = -3199.293799999999
$ws.Range("N126").Value = -42940.001
